{"js": "// Remove the \"Table of Contents\" building block (the w:sdt wrapping the\n// TOC heading paragraph + the TOC field code) that sits at the very start\n// of the document body.\n//\n// In the Word/Office.js object model this block is a content control of\n// subtype \"BuildingBlockGallery\" (docPartGallery = \"Table of Contents\").\n// `document.contentControls` only surfaces rich-text content controls, so\n// we use `document.getContentControls()` which enumerates *all* SDT types\n// (matches the Word desktop behaviour for `Document.ContentControls`).\n\nconst contentControls = context.document.getContentControls();\ncontentControls.load(\"items\");\nawait context.sync();\n\ncontentControls.items.forEach((cc) => {\n  cc.load(\"subtype\");\n});\nawait context.sync();\n\ncontentControls.items.forEach((cc) => {\n  if (cc.subtype === \"BuildingBlockGallery\") {\n    // keepContent = false: delete the control AND its contents (the TOC\n    // heading paragraph and the TOC field paragraph), exactly removing the\n    // whole <w:sdt>...</w:sdt> node from document.xml.\n    cc.delete(false);\n  }\n});\nawait context.sync();\n", "ps1": "# Remove the \"Table of Contents\" building block (the <w:sdt> wrapping the\n# TOC heading paragraph + the TOC field code) that sits at the very start\n# of the document body.\n#\n# In Word's object model this block is a content control of type\n# wdContentControlBuildingBlockGallery (5) whose gallery is \"Table of\n# Contents\" (it is inserted by the built-in \"Table of Contents\" quick-part /\n# the References > Table of Contents gallery). We find it through\n# $d.ContentControls and delete it together with its content.\n\n$d = $word.ActiveDocument\n\n# Walk backwards so deleting one control does not disturb the indices of\n# the ones still to be visited.\nfor ($i = $d.ContentControls.Count; $i -ge 1; $i--) {\n    $cc = $d.ContentControls.Item($i)\n\n    if ($cc.Type -eq 5) {\n        # wdContentControlBuildingBlockGallery\n        if ($cc.LockContentControl) {\n            $cc.LockContentControl = $false\n        }\n        # DeleteContents = $true removes the control AND everything inside\n        # it (the \"Table of Contents\" heading paragraph and the TOC field\n        # paragraph), i.e. the whole <w:sdt> node disappears from the XML.\n        $cc.Delete($true)\n    }\n}\n"}
